$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 600
$wsExhibit.Range("F7").Value = 1802
$wsExhibit.Range("F8").Value = 98

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 107

# Sheet "全部类型" (All types - aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 600
$wsAll.Range("F7").Value = 107
$wsAll.Range("F11").Value = 1802
$wsAll.Range("F12").Value = 98

$wb.Save()
